$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 90
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6481.25
$ws.Range("I70").Value = 5112
$ws.Range("J70").Value = 7850.5
$ws.Range("K70").Value = 15336
$ws.Range("L70").Value = 23551.5
$ws.Range("M70").Value = -15066
$ws.Range("N70").Value = -24091.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6481.25
$ws.Range("I73").Value = 5112
$ws.Range("J73").Value = 7850.5
$ws.Range("K73").Value = 15336
$ws.Range("L73").Value = 23551.5
$ws.Range("M73").Value = -14400
$ws.Range("N73").Value = -25423.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 512.03705
$ws.Range("I92").Value = 480.9091
$ws.Range("J92").Value = 649
$ws.Range("K92").Value = 480.9091
$ws.Range("L92").Value = 649
$ws.Range("M92").Value = 767.0908999999999
$ws.Range("N92").Value = -3145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4985.3
$ws.Range("I100").Value = 5675.5
$ws.Range("J100").Value = 3950
$ws.Range("K100").Value = 5675.5
$ws.Range("L100").Value = 3950
$ws.Range("M100").Value = -5134.5
$ws.Range("N100").Value = -5032

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3103.4211
$ws.Range("I138").Value = 996
$ws.Range("J138").Value = 3220.5
$ws.Range("K138").Value = 2988
$ws.Range("L138").Value = 9661.5
$ws.Range("M138").Value = 2152
$ws.Range("N138").Value = -19941.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 11599.2
$ws.Range("I141").Value = 14665.667
$ws.Range("J141").Value = 6999.5
$ws.Range("K141").Value = 43997.001
$ws.Range("L141").Value = 20998.5
$ws.Range("M141").Value = -38817.001
$ws.Range("N141").Value = -31358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 26376832
$ws.Range("I61").Value = 62503870
$ws.Range("J61").Value = 102622.37
$ws.Range("K61").Value = 62503870
$ws.Range("L61").Value = 102622.37
$ws.Range("M61").Value = -62503658
$ws.Range("N61").Value = -103046.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2870.3333
$ws.Range("I110").Value = 2555.5
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 2555.5
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = -510.5
$ws.Range("N110").Value = -7590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1252.7273
$ws.Range("I122").Value = 1168
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 3504
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -1054
$ws.Range("N122").Value = -11200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3652.8948
$ws.Range("I132").Value = 2436.2
$ws.Range("J132").Value = 6515.706
$ws.Range("K132").Value = 7308.599999999999
$ws.Range("L132").Value = 19547.118
$ws.Range("M132").Value = -4778.599999999999
$ws.Range("N132").Value = -24607.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 64998.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 64998.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 64998.5
$ws.Range("N133").Value = -70058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 26376832
$ws.Range("I136").Value = 62503870
$ws.Range("J136").Value = 102622.37
$ws.Range("K136").Value = 187511610
$ws.Range("L136").Value = 307867.11
$ws.Range("M136").Value = -187509060
$ws.Range("N136").Value = -312967.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4369.5
$ws.Range("I5").Value = 5676.3335
$ws.Range("J5").Value = 449
$ws.Range("K5").Value = 5676.3335
$ws.Range("L5").Value = 449
$ws.Range("M5").Value = -5563.3335
$ws.Range("N5").Value = -675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2828.4666
$ws.Range("I86").Value = 1547.7
$ws.Range("J86").Value = 5390
$ws.Range("K86").Value = 1547.7
$ws.Range("L86").Value = 5390
$ws.Range("M86").Value = -424.7
$ws.Range("N86").Value = -7636

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2828.4666
$ws.Range("I89").Value = 1547.7
$ws.Range("J89").Value = 5390
$ws.Range("K89").Value = 7738.5
$ws.Range("L89").Value = 26950
$ws.Range("M89").Value = -2122.5
$ws.Range("N89").Value = -38182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2298.3809
$ws.Range("I99").Value = 1864.7778
$ws.Range("J99").Value = 4900
$ws.Range("K99").Value = 1864.7778
$ws.Range("L99").Value = 4900
$ws.Range("M99").Value = -366.7778000000001
$ws.Range("N99").Value = -7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 133802.62
$ws.Range("I134").Value = 5139.6665
$ws.Range("J134").Value = 211000.4
$ws.Range("K134").Value = 15418.9995
$ws.Range("L134").Value = 633001.2
$ws.Range("M134").Value = -12883.9995
$ws.Range("N134").Value = -638071.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 18400
$ws.Range("I74").Value = 16000
$ws.Range("J74").Value = 19000
$ws.Range("K74").Value = 16000
$ws.Range("L74").Value = 19000
$ws.Range("M74").Value = -15126
$ws.Range("N74").Value = -20748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 18400
$ws.Range("I77").Value = 16000
$ws.Range("J77").Value = 19000
$ws.Range("K77").Value = 48000
$ws.Range("L77").Value = 57000
$ws.Range("M77").Value = -43632
$ws.Range("N77").Value = -65736

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5767.067
$ws.Range("I137").Value = 6361
$ws.Range("J137").Value = 5247.375
$ws.Range("K137").Value = 19083
$ws.Range("L137").Value = 15742.125
$ws.Range("M137").Value = -13983
$ws.Range("N137").Value = -25942.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4233.7
$ws.Range("I102").Value = 3148.6667
$ws.Range("J102").Value = 13999
$ws.Range("K102").Value = 3148.6667
$ws.Range("L102").Value = 13999
$ws.Range("M102").Value = -1526.6667
$ws.Range("N102").Value = -17243

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2041.625
$ws.Range("I107").Value = 2041.625
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2041.625
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -121.625
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4246
$ws.Range("I113").Value = 3787.1428
$ws.Range("J113").Value = 5316.6665
$ws.Range("K113").Value = 3787.1428
$ws.Range("L113").Value = 5316.6665
$ws.Range("M113").Value = -1617.1428
$ws.Range("N113").Value = -9656.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 50002892
$ws.Range("I132").Value = 66669590
$ws.Range("J132").Value = 2792.6
$ws.Range("K132").Value = 200008770
$ws.Range("L132").Value = 8377.799999999999
$ws.Range("M132").Value = -200006240
$ws.Range("N132").Value = -13437.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1731.25
$ws.Range("I16").Value = 1528.6154
$ws.Range("J16").Value = 1970.7273
$ws.Range("K16").Value = 1528.6154
$ws.Range("L16").Value = 1970.7273
$ws.Range("M16").Value = -1358.6154
$ws.Range("N16").Value = -2310.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2883
$ws.Range("I40").Value = 2036.4
$ws.Range("J40").Value = 4999.5
$ws.Range("K40").Value = 2036.4
$ws.Range("L40").Value = 4999.5
$ws.Range("M40").Value = -1900.4
$ws.Range("N40").Value = -5271.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 408955.28
$ws.Range("I132").Value = 386222.62
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 1158667.86
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -1156137.86
$ws.Range("N132").Value = -3005075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 109247.18
$ws.Range("I136").Value = 6452.5
$ws.Range("J136").Value = 232600.8
$ws.Range("K136").Value = 19357.5
$ws.Range("L136").Value = 697802.3999999999
$ws.Range("M136").Value = -16807.5
$ws.Range("N136").Value = -702902.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

Write-Host "edit complete"
